# HG0B1CEU BOM — "Correcting some documentaion issues"
#
# The Manufacturer (col C) / Manufacturer Part Number (col D) values were
# transposed for the Harwin connectors (rows 9 & 10) and the MCC diode
# (row 23). Swap each pair back into the correct column.
#
# D11's part number also gets re-formatted to match the plain
# center+vertically-centered style used elsewhere in that column (e.g.
# C12/D12) instead of the custom numeric format it had before.
#
# Finally, restore the on-screen view: scrolled so row 4 is the top row,
# with C23 as the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-CellValues {
    param(
        [string]$Address1,
        [string]$Address2
    )
    $range1 = $ws.Range($Address1)
    $range2 = $ws.Range($Address2)
    $value1 = $range1.Value2
    $value2 = $range2.Value2
    $range1.Value2 = $value2
    $range2.Value2 = $value1
}

# Fix Manufacturer / Manufacturer Part Number column swaps
Swap-CellValues "C9" "D9"
Swap-CellValues "C10" "D10"
Swap-CellValues "C23" "D23"

# D11 (Manufacturer Part Number for the USB connector) should use the same
# formatting as the rest of that column (center + vertical-center, general
# number format) rather than its previous custom numeric format.
$ws.Range("C12").Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Restore the saved scroll position / selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C23").Select()

Write-Output "HG0B1CEU BOM documentation corrections applied"
